$d = $word.ActiveDocument

function Replace-Exact($findText, $replaceText) {
    $rng = $d.Content
    $found = $rng.Find.Execute(
        $findText,
        $true, $false, $false, $false, $false,
        $true, 1, $false
    )
    if ($found) {
        $rng.Text = $replaceText
    }
}

Replace-Exact '"Sans auth, accès refusé."' '"Sans authentification, l''accès est refusé."'

Replace-Exact 'src/services/auth_service.py (lignes 33-37 + 135-143)' 'src/services/token_service.py (lignes 31-32 + 73-85)'

Replace-Exact 'src/cli/commands/user_commands.py (ligne ~25)' 'src/cli/commands/user_commands.py (lignes 13-15)'

Replace-Exact 'src/models/user.py (lignes 56-60)' 'src/services/password_hashing_service.py (lignes 38-41)'

Replace-Exact 'src/cli/commands/client_commands.py (lignes 72-79)' 'src/cli/commands/client_commands.py (lignes 76-78)'
